# "Set up CPR filter test cases"
#
# 1) Rename "EP2040 Memory Map" -> "RP2040 Memory Map" (typo fix).
# 2) Move the selection on that sheet to E15 (it no longer holds the tab focus).
# 3) Add a new worksheet "Mode S Beast RSSI Byte" after the memory-map sheet
#    and populate it with the RSSI/dBm test-case table.

$wb = $excel.ActiveWorkbook

# --- Fix sheet name -------------------------------------------------------
$memMap = $wb.Worksheets.Item("EP2040 Memory Map")
$memMap.Name = "RP2040 Memory Map"
$memMap.Range("E15").Select() | Out-Null

# --- New sheet with CPR / Mode S Beast RSSI byte test cases ---------------
$ws = $wb.Worksheets.Add($null, $memMap)
$ws.Name = "Mode S Beast RSSI Byte"

$ws.Range("A1").Value = "Parameter"
$ws.Range("B1").Value = "Value"
$ws.Range("C1").Value = "Unit"
$ws.Range("D1").Value = "Note"

$ws.Range("A2").Value = "p_min"
$ws.Range("B2").Value = -95
$ws.Range("C2").Value = "dBm"

$ws.Range("A3").Value = "p_max"
$ws.Range("B3").Value = -45
$ws.Range("C3").Value = "dBm"

$ws.Range("A5").Value = "rssi_min"
$ws.Range("B5").Formula = "=255*10^((B2+45)/10)"

$ws.Range("A6").Value = "rssi_max"
$ws.Range("B6").Formula = "=255*10^((B3+45)/10)"

$ws.Range("A8").Value = "p_mid"
$ws.Range("B8").Value = -70

$ws.Range("A9").Value = "rssi_mid"
$ws.Range("B9").Formula = "=255*10^((B8+45)/10)"

$ws.Columns.Item(1).ColumnWidth = 24.75

$ws.Range("B9").Select() | Out-Null
